# Insert a new "Completed Tasks" slide (with a summary table) right
# before the existing "Process model" slide. PowerPoint re-numbers the
# slides that follow automatically; the new slide receives the next
# available sldId.

$p = $ppt.ActivePresentation

# --- New slide 7: "Completed Tasks" -----------------------------------
# Layout 11 = ppLayoutTitleOnly, matching the other content slides in
# this deck (slideLayout6.xml / "Title Only").
$s = $p.Slides.Add(7, 11)

$titleRange = $s.Shapes.Title.TextFrame.TextRange
$titleRange.Text = "Completed Tasks "
$titleRange.ParagraphFormat.Alignment = 2  # ppAlignCenter

# --- Table -------------------------------------------------------------
# EMU -> point conversion (1 pt = 12700 EMU) since Shapes.AddTable /
# Width / Height take points.
$tbl = $s.Shapes.AddTable(9, 2, 2967553 / 12700, 2152175 / 12700, 5412967 / 12700, 4110552 / 12700)

$tbl.Table.ApplyStyle("{616DA210-FB5B-4158-B5E0-FEB733F419BA}")
$tbl.Table.FirstRow = $false
$tbl.Table.HorizBanding = $false

$tbl.Table.Columns.Item(1).Width = 1870700 / 12700
$tbl.Table.Columns.Item(2).Width = 3542267 / 12700

$rows = @(
    @("Task", "Who Completed"),
    @("GitHub Setup", "Josh"),
    @("Use Cases/Requirements", "Christian"),
    @("Gantt Chart ", "Travis"),
    @("SPMP ", "Travis, Christian"),
    @("Technical Status Presentation", "Travis"),
    @("Integration/Unit Testing Skeleton", " "),
    @("Website/Game Development ", "Josh, Travis"),
    @(" ", " ")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowIndex = $i + 1
    $left = $rows[$i][0]
    $right = $rows[$i][1]

    $leftCell = $tbl.Table.Cell($rowIndex, 1).Shape.TextFrame.TextRange
    $leftCell.Text = $left
    $leftCell.Font.Size = 12

    $rightCell = $tbl.Table.Cell($rowIndex, 2).Shape.TextFrame.TextRange
    $rightCell.Text = $right
    $rightCell.Font.Size = 12
    if ($rowIndex -ge 2 -and $rowIndex -le 8) {
        $rightCell.ParagraphFormat.Alignment = 2  # ppAlignCenter
    }
}
